$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '79.046.83'
Set-TextValue 'E2' '  +3.21%  '
Set-TextValue 'D3' '3.181.22'
Set-TextValue 'E3' '  +5.05%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '204.99'
Set-TextValue 'E5' '  +1.74%  '
Set-TextValue 'D6' '630.68'
Set-TextValue 'E6' '  +0.17%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'E8' '  +7.55%  '
Set-TextValue 'D9' '0.584'
Set-TextValue 'E9' '  +5.44%  '
Set-TextValue 'D10' '3.177.94'
Set-TextValue 'E10' '  +4.92%  '
Set-TextValue 'D11' '0.578'
Set-TextValue 'E11' '  +32.29%  '
Set-TextValue 'E12' '  +2.81%  '
Set-TextValue 'E13' '  +5.15%  '
Set-TextValue 'D14' '3.767.04'
Set-TextValue 'E14' '  +5.06%  '
Set-TextValue 'D15' '0.0000226'
Set-TextValue 'E15' '  +15.30%  '
Set-TextValue 'D16' '31.54'
Set-TextValue 'E16' '  +7.19%  '
Set-TextValue 'D17' '78.870.94'
Set-TextValue 'E17' '  +3.11%  '
Set-TextValue 'D18' '3.178.72'
Set-TextValue 'E18' '  +5.48%  '
Set-TextValue 'D19' '14.47'
Set-TextValue 'E19' '  +7.03%  '
Set-TextValue 'B20' 'SuiNetwork'
Set-TextValue 'C20' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D20' '2.97'
Set-TextValue 'E20' '  +30.00%  '
Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '9.24'
Set-TextValue 'E21' '  +2.57%  '
Set-TextValue 'D22' '427.23'
Set-TextValue 'E22' '  +13.55%  '
Set-TextValue 'D23' '4.95'
Set-TextValue 'E23' '  +13.38%  '
Set-TextValue 'D24' '6.85'
Set-TextValue 'E24' '  +6.01%  '
Set-TextValue 'D25' '3.350.94'
Set-TextValue 'E25' '  +5.20%  '
Set-TextValue 'D26' '4.77'
Set-TextValue 'E26' '  +8.45%  '
Set-TextValue 'D27' '11.06'
Set-TextValue 'E27' '  +10.97%  '
Set-TextValue 'D28' '76.34'
Set-TextValue 'E28' '  +4.10%  '
Set-TextValue 'E29' '  +0.41%  '
Set-TextValue 'E30' '  +2.82%  '
Set-TextValue 'D31' '0.999'
Set-TextValue 'E31' '  +0.08%  '
Set-TextValue 'D32' '8.97'
Set-TextValue 'E32' '  +7.40%  '
Set-TextValue 'E33' '  +4.02%  '
Set-TextValue 'D34' '520.99'
Set-TextValue 'E34' '  +2.56%  '
Set-TextValue 'E35' '  +2.15%  '
Set-TextValue 'D36' '22.97'
Set-TextValue 'E36' '  +10.85%  '
Set-TextValue 'B37' 'Kaspa'
Set-TextValue 'C37' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D37' '0.136'
Set-TextValue 'E37' '  +20.35%  '
Set-TextValue 'B38' 'Cronos'
Set-TextValue 'C38' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D38' '0.125'
Set-TextValue 'E38' '  +18.73%  '
Set-TextValue 'E39' '  -0.10%  '
Set-TextValue 'D40' '0.400'
Set-TextValue 'E40' '  +3.76%  '
Set-TextValue 'D41' '164.05'
Set-TextValue 'E41' '  +0.03%  '
Set-TextValue 'D42' '20.00'
Set-TextValue 'E42' '  -0.09%  '
Set-TextValue 'B43' 'USDe'
Set-TextValue 'C43' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D43' '1.00'
Set-TextValue 'E43' '  +0.01%  '
Set-TextValue 'B44' 'Aave'
Set-TextValue 'C44' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D44' '191.91'
Set-TextValue 'E44' '  +2.10%  '
Set-TextValue 'E45' '  +5.20%  '
Set-TextValue 'D46' '0.817'
Set-TextValue 'E46' '  +14.24%  '
Set-TextValue 'D47' '1.79'
Set-TextValue 'E47' '  +6.20%  '
Set-TextValue 'E48' '  +3.86%  '
Set-TextValue 'D49' '42.55'
Set-TextValue 'E49' '  +0.12%  '
Set-TextValue 'D50' '2.51'
Set-TextValue 'E50' '  +3.43%  '
Set-TextValue 'B51' 'InjectiveProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D51' '25.34'
Set-TextValue 'E51' '  +11.83%  '
